$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 11:11:18"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 11:11:13"
$wsZhCn.Range("K2").Value = "2016-09-04 11:11:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 11:11:18"
$wsDeDe.Range("K2").Value = "2016-09-04 11:11:38"
